# Edit script: updates lockdown-effectiveness weighting (Offices column Z
# excluded from the weight total) and appends 12 new date rows
# (9/30/2020 .. 10/11/2020) copied from the last existing row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Converted Data")

# --- Step 1: Offices (Z) no longer counted in the weight sum ---
$ws.Range("Z7").Value = 0
$ws.Range("AE7").Value = 12

# --- Step 2: recalculate LockdownEffectiveness (AE) for existing rows ---
$ws.Range("AE20").Value = 0.08333333333333333
$ws.Range("AE21").Value = 0.08333333333333333
$ws.Range("AE22").Value = 0.08333333333333333
$ws.Range("AE23").Value = 0.08333333333333333
$ws.Range("AE24").Value = 0.30555555554166663
$ws.Range("AE25").Value = 0.388888888875
$ws.Range("AE26").Value = 0.388888888875
$ws.Range("AE27").Value = 0.5694444444250001
$ws.Range("AE28").Value = 0.6527777777583333
$ws.Range("AE29").Value = 0.6527777777583333
$ws.Range("AE30").Value = 0.6527777777583333
$ws.Range("AE31").Value = 0.8333333333083334
$ws.Range("AE32").Value = 0.8333333333083334
$ws.Range("AE33").Value = 0.8333333333083334
$ws.Range("AE34").Value = 0.8333333333083334
$ws.Range("AE35").Value = 0.8333333333083334
$ws.Range("AE36").Value = 0.8333333333083334
$ws.Range("AE37").Value = 0.8333333333083334
$ws.Range("AE38").Value = 0.8333333333083334
$ws.Range("AE39").Value = 0.8333333333083334
$ws.Range("AE40").Value = 0.8333333333083334
$ws.Range("AE41").Value = 0.8333333333083334
$ws.Range("AE42").Value = 0.8333333333083334
$ws.Range("AE43").Value = 0.8333333333083334
$ws.Range("AE44").Value = 0.8333333333083334
$ws.Range("AE45").Value = 0.8333333333083334
$ws.Range("AE46").Value = 0.8333333333083334
$ws.Range("AE47").Value = 0.8333333333083334
$ws.Range("AE48").Value = 0.8333333333083334
$ws.Range("AE49").Value = 0.8333333333083334
$ws.Range("AE50").Value = 0.8333333333083334
$ws.Range("AE51").Value = 0.8333333333083334
$ws.Range("AE52").Value = 0.8333333333083334
$ws.Range("AE53").Value = 0.8333333333083334
$ws.Range("AE54").Value = 0.8333333333083334
$ws.Range("AE55").Value = 0.8333333333083334
$ws.Range("AE56").Value = 0.8333333333083334
$ws.Range("AE57").Value = 0.8333333333083334
$ws.Range("AE58").Value = 0.8333333333083334
$ws.Range("AE59").Value = 0.9166666666416666
$ws.Range("AE60").Value = 0.9166666666416666
$ws.Range("AE61").Value = 0.9166666666416666
$ws.Range("AE62").Value = 0.9166666666416666
$ws.Range("AE63").Value = 0.9166666666416666
$ws.Range("AE64").Value = 0.9166666666416666
$ws.Range("AE65").Value = 0.9166666666416666
$ws.Range("AE66").Value = 0.9166666666416666
$ws.Range("AE67").Value = 0.9166666666416666
$ws.Range("AE68").Value = 0.9166666666416666
$ws.Range("AE69").Value = 0.9166666666416666
$ws.Range("AE70").Value = 0.9166666666416666
$ws.Range("AE71").Value = 0.9166666666416666
$ws.Range("AE72").Value = 0.9166666666416666
$ws.Range("AE73").Value = 0.9166666666416666
$ws.Range("AE74").Value = 0.9166666666416666
$ws.Range("AE75").Value = 0.9166666666416666
$ws.Range("AE76").Value = 0.9166666666416666
$ws.Range("AE77").Value = 0.9166666666416666
$ws.Range("AE78").Value = 0.9166666666416666
$ws.Range("AE79").Value = 0.9166666666416666
$ws.Range("AE80").Value = 0.9166666666416666
$ws.Range("AE81").Value = 0.9166666666416666
$ws.Range("AE82").Value = 0.9166666666416666
$ws.Range("AE83").Value = 0.9166666666416666
$ws.Range("AE84").Value = 0.9166666666416666
$ws.Range("AE85").Value = 0.9166666666416666
$ws.Range("AE86").Value = 0.9166666666416666
$ws.Range("AE87").Value = 0.9166666666416666
$ws.Range("AE88").Value = 0.9166666666416666
$ws.Range("AE89").Value = 0.9166666666416666
$ws.Range("AE90").Value = 0.6666666666416666
$ws.Range("AE91").Value = 0.6666666666416666
$ws.Range("AE92").Value = 0.6388888888666667
$ws.Range("AE93").Value = 0.6388888888666667
$ws.Range("AE94").Value = 0.6388888888666667
$ws.Range("AE95").Value = 0.6388888888666667
$ws.Range("AE96").Value = 0.6388888888666667
$ws.Range("AE97").Value = 0.6388888888666667
$ws.Range("AE98").Value = 0.6388888888666667
$ws.Range("AE99").Value = 0.6388888888666667
$ws.Range("AE100").Value = 0.6111111110916666
$ws.Range("AE101").Value = 0.6111111110916666
$ws.Range("AE102").Value = 0.5777777777583334
$ws.Range("AE103").Value = 0.5777777777583334
$ws.Range("AE104").Value = 0.5777777777583334
$ws.Range("AE105").Value = 0.5777777777583334
$ws.Range("AE106").Value = 0.5777777777583334
$ws.Range("AE107").Value = 0.5777777777583334
$ws.Range("AE108").Value = 0.5777777777583334
$ws.Range("AE109").Value = 0.5777777777583334
$ws.Range("AE110").Value = 0.5777777777583334
$ws.Range("AE111").Value = 0.5777777777583334
$ws.Range("AE112").Value = 0.5777777777583334
$ws.Range("AE113").Value = 0.5777777777583334
$ws.Range("AE114").Value = 0.5777777777583334
$ws.Range("AE115").Value = 0.5777777777583334
$ws.Range("AE116").Value = 0.5777777777583334
$ws.Range("AE117").Value = 0.5777777777583334
$ws.Range("AE118").Value = 0.30555555554999997
$ws.Range("AE119").Value = 0.30555555554999997
$ws.Range("AE120").Value = 0.30555555554999997
$ws.Range("AE121").Value = 0.30555555554999997
$ws.Range("AE122").Value = 0.30555555554999997
$ws.Range("AE123").Value = 0.30555555554999997
$ws.Range("AE124").Value = 0.30555555554999997
$ws.Range("AE125").Value = 0.30555555554999997
$ws.Range("AE126").Value = 0.30555555554999997
$ws.Range("AE127").Value = 0.30555555554999997
$ws.Range("AE128").Value = 0.30555555554999997
$ws.Range("AE129").Value = 0.30555555554999997
$ws.Range("AE130").Value = 0.30555555554999997
$ws.Range("AE131").Value = 0.30555555554999997
$ws.Range("AE132").Value = 0.30555555554999997
$ws.Range("AE133").Value = 0.30555555554999997
$ws.Range("AE134").Value = 0.30555555554999997
$ws.Range("AE135").Value = 0.30555555554999997
$ws.Range("AE136").Value = 0.30555555554999997
$ws.Range("AE137").Value = 0.30555555554999997
$ws.Range("AE138").Value = 0.30555555554999997
$ws.Range("AE139").Value = 0.277777777775
$ws.Range("AE140").Value = 0.277777777775
$ws.Range("AE141").Value = 0.277777777775
$ws.Range("AE142").Value = 0.277777777775
$ws.Range("AE143").Value = 0.277777777775
$ws.Range("AE144").Value = 0.277777777775
$ws.Range("AE145").Value = 0.277777777775
$ws.Range("AE146").Value = 0.277777777775
$ws.Range("AE147").Value = 0.277777777775
$ws.Range("AE148").Value = 0.277777777775
$ws.Range("AE149").Value = 0.277777777775
$ws.Range("AE150").Value = 0.277777777775
$ws.Range("AE151").Value = 0.277777777775
$ws.Range("AE152").Value = 0.277777777775
$ws.Range("AE153").Value = 0.277777777775
$ws.Range("AE154").Value = 0.277777777775
$ws.Range("AE155").Value = 0.277777777775
$ws.Range("AE156").Value = 0.277777777775
$ws.Range("AE157").Value = 0.277777777775
$ws.Range("AE158").Value = 0.277777777775
$ws.Range("AE159").Value = 0.277777777775
$ws.Range("AE160").Value = 0.277777777775
$ws.Range("AE161").Value = 0.277777777775
$ws.Range("AE162").Value = 0.277777777775
$ws.Range("AE163").Value = 0.277777777775
$ws.Range("AE164").Value = 0.277777777775
$ws.Range("AE165").Value = 0.277777777775
$ws.Range("AE166").Value = 0.277777777775
$ws.Range("AE167").Value = 0.277777777775
$ws.Range("AE168").Value = 0.277777777775
$ws.Range("AE169").Value = 0.277777777775
$ws.Range("AE170").Value = 0.277777777775
$ws.Range("AE171").Value = 0.277777777775
$ws.Range("AE172").Value = 0.277777777775
$ws.Range("AE173").Value = 0.277777777775
$ws.Range("AE174").Value = 0.277777777775
$ws.Range("AE175").Value = 0.277777777775
$ws.Range("AE176").Value = 0.277777777775
$ws.Range("AE177").Value = 0.277777777775
$ws.Range("AE178").Value = 0.277777777775
$ws.Range("AE179").Value = 0.277777777775
$ws.Range("AE180").Value = 0.277777777775
$ws.Range("AE181").Value = 0.277777777775
$ws.Range("AE182").Value = 0.277777777775
$ws.Range("AE183").Value = 0.277777777775
$ws.Range("AE184").Value = 0.277777777775
$ws.Range("AE185").Value = 0.277777777775
$ws.Range("AE186").Value = 0.277777777775
$ws.Range("AE187").Value = 0.277777777775
$ws.Range("AE188").Value = 0.277777777775
$ws.Range("AE189").Value = 0.277777777775
$ws.Range("AE190").Value = 0.277777777775
$ws.Range("AE191").Value = 0.277777777775
$ws.Range("AE192").Value = 0.277777777775
$ws.Range("AE193").Value = 0.277777777775
$ws.Range("AE194").Value = 0.277777777775
$ws.Range("AE195").Value = 0.277777777775
$ws.Range("AE196").Value = 0.277777777775
$ws.Range("AE197").Value = 0.277777777775
$ws.Range("AE198").Value = 0.277777777775
$ws.Range("AE199").Value = 0.277777777775
$ws.Range("AE200").Value = 0.277777777775
$ws.Range("AE201").Value = 0.277777777775
$ws.Range("AE202").Value = 0.277777777775
$ws.Range("AE203").Value = 0.277777777775
$ws.Range("AE204").Value = 0.277777777775
$ws.Range("AE205").Value = 0.277777777775
$ws.Range("AE206").Value = 0.277777777775
$ws.Range("AE207").Value = 0.277777777775
$ws.Range("AE208").Value = 0.277777777775
$ws.Range("AE209").Value = 0.277777777775
$ws.Range("AE210").Value = 0.277777777775
$ws.Range("AE211").Value = 0.277777777775
$ws.Range("AE212").Value = 0.277777777775
$ws.Range("AE213").Value = 0.277777777775
$ws.Range("AE214").Value = 0.277777777775
$ws.Range("AE215").Value = 0.277777777775
$ws.Range("AE216").Value = 0.277777777775
$ws.Range("AE217").Value = 0.277777777775
$ws.Range("AE218").Value = 0.277777777775
$ws.Range("AE219").Value = 0.277777777775
$ws.Range("AE220").Value = 0.277777777775
$ws.Range("AE221").Value = 0.277777777775

# --- Step 3: append 12 new date rows (222-233) ---
$ws.Range("A222").Value = "'9/30/2020"
$ws.Range("A223").Value = "'10/1/2020"
$ws.Range("A224").Value = "'10/2/2020"
$ws.Range("A225").Value = "'10/3/2020"
$ws.Range("A226").Value = "'10/4/2020"
$ws.Range("A227").Value = "'10/5/2020"
$ws.Range("A228").Value = "'10/6/2020"
$ws.Range("A229").Value = "'10/7/2020"
$ws.Range("A230").Value = "'10/8/2020"
$ws.Range("A231").Value = "'10/9/2020"
$ws.Range("A232").Value = "'10/10/2020"
$ws.Range("A233").Value = "'10/11/2020"

# Re-apply the header/date-column cell format (border + bold + centered)
# from the last existing date row so the new cells match it exactly.
$ws.Range("A221").Copy()
$ws.Range("A222:A233").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Step 3b: flag columns (B..AD) -- same pattern as the last existing row
$ws.Range("B222").Value = 0
$ws.Range("C222").Value = 0
$ws.Range("D222").Value = 1
$ws.Range("E222").Value = 1
$ws.Range("F222").Value = 1
$ws.Range("G222").Value = 0
$ws.Range("H222").Value = 0
$ws.Range("I222").Value = 0
$ws.Range("J222").Value = 0
$ws.Range("K222").Value = 0
$ws.Range("L222").Value = 0
$ws.Range("M222").Value = 0
$ws.Range("N222").Value = 0
$ws.Range("O222").Value = 0
$ws.Range("P222").Value = 0
$ws.Range("Q222").Value = 0
$ws.Range("R222").Value = 0
$ws.Range("S222").Value = 0
$ws.Range("T222").Value = 0
$ws.Range("U222").Value = 0
$ws.Range("V222").Value = 0
$ws.Range("W222").Value = 0
$ws.Range("X222").Value = 0
$ws.Range("Y222").Value = 0
$ws.Range("Z222").Value = 0
$ws.Range("AA222").Value = 0
$ws.Range("AB222").Value = 1
$ws.Range("AC222").Value = 1
$ws.Range("AD222").Value = 0
$ws.Range("B223").Value = 0
$ws.Range("C223").Value = 0
$ws.Range("D223").Value = 1
$ws.Range("E223").Value = 1
$ws.Range("F223").Value = 1
$ws.Range("G223").Value = 0
$ws.Range("H223").Value = 0
$ws.Range("I223").Value = 0
$ws.Range("J223").Value = 0
$ws.Range("K223").Value = 0
$ws.Range("L223").Value = 0
$ws.Range("M223").Value = 0
$ws.Range("N223").Value = 0
$ws.Range("O223").Value = 0
$ws.Range("P223").Value = 0
$ws.Range("Q223").Value = 0
$ws.Range("R223").Value = 0
$ws.Range("S223").Value = 0
$ws.Range("T223").Value = 0
$ws.Range("U223").Value = 0
$ws.Range("V223").Value = 0
$ws.Range("W223").Value = 0
$ws.Range("X223").Value = 0
$ws.Range("Y223").Value = 0
$ws.Range("Z223").Value = 0
$ws.Range("AA223").Value = 0
$ws.Range("AB223").Value = 1
$ws.Range("AC223").Value = 1
$ws.Range("AD223").Value = 0
$ws.Range("B224").Value = 0
$ws.Range("C224").Value = 0
$ws.Range("D224").Value = 1
$ws.Range("E224").Value = 1
$ws.Range("F224").Value = 1
$ws.Range("G224").Value = 0
$ws.Range("H224").Value = 0
$ws.Range("I224").Value = 0
$ws.Range("J224").Value = 0
$ws.Range("K224").Value = 0
$ws.Range("L224").Value = 0
$ws.Range("M224").Value = 0
$ws.Range("N224").Value = 0
$ws.Range("O224").Value = 0
$ws.Range("P224").Value = 0
$ws.Range("Q224").Value = 0
$ws.Range("R224").Value = 0
$ws.Range("S224").Value = 0
$ws.Range("T224").Value = 0
$ws.Range("U224").Value = 0
$ws.Range("V224").Value = 0
$ws.Range("W224").Value = 0
$ws.Range("X224").Value = 0
$ws.Range("Y224").Value = 0
$ws.Range("Z224").Value = 0
$ws.Range("AA224").Value = 0
$ws.Range("AB224").Value = 1
$ws.Range("AC224").Value = 1
$ws.Range("AD224").Value = 0
$ws.Range("B225").Value = 0
$ws.Range("C225").Value = 0
$ws.Range("D225").Value = 1
$ws.Range("E225").Value = 1
$ws.Range("F225").Value = 1
$ws.Range("G225").Value = 0
$ws.Range("H225").Value = 0
$ws.Range("I225").Value = 0
$ws.Range("J225").Value = 0
$ws.Range("K225").Value = 0
$ws.Range("L225").Value = 0
$ws.Range("M225").Value = 0
$ws.Range("N225").Value = 0
$ws.Range("O225").Value = 0
$ws.Range("P225").Value = 0
$ws.Range("Q225").Value = 0
$ws.Range("R225").Value = 0
$ws.Range("S225").Value = 0
$ws.Range("T225").Value = 0
$ws.Range("U225").Value = 0
$ws.Range("V225").Value = 0
$ws.Range("W225").Value = 0
$ws.Range("X225").Value = 0
$ws.Range("Y225").Value = 0
$ws.Range("Z225").Value = 0
$ws.Range("AA225").Value = 0
$ws.Range("AB225").Value = 1
$ws.Range("AC225").Value = 1
$ws.Range("AD225").Value = 0
$ws.Range("B226").Value = 0
$ws.Range("C226").Value = 0
$ws.Range("D226").Value = 1
$ws.Range("E226").Value = 1
$ws.Range("F226").Value = 1
$ws.Range("G226").Value = 0
$ws.Range("H226").Value = 0
$ws.Range("I226").Value = 0
$ws.Range("J226").Value = 0
$ws.Range("K226").Value = 0
$ws.Range("L226").Value = 0
$ws.Range("M226").Value = 0
$ws.Range("N226").Value = 0
$ws.Range("O226").Value = 0
$ws.Range("P226").Value = 0
$ws.Range("Q226").Value = 0
$ws.Range("R226").Value = 0
$ws.Range("S226").Value = 0
$ws.Range("T226").Value = 0
$ws.Range("U226").Value = 0
$ws.Range("V226").Value = 0
$ws.Range("W226").Value = 0
$ws.Range("X226").Value = 0
$ws.Range("Y226").Value = 0
$ws.Range("Z226").Value = 0
$ws.Range("AA226").Value = 0
$ws.Range("AB226").Value = 1
$ws.Range("AC226").Value = 1
$ws.Range("AD226").Value = 0
$ws.Range("B227").Value = 0
$ws.Range("C227").Value = 0
$ws.Range("D227").Value = 1
$ws.Range("E227").Value = 1
$ws.Range("F227").Value = 1
$ws.Range("G227").Value = 0
$ws.Range("H227").Value = 0
$ws.Range("I227").Value = 0
$ws.Range("J227").Value = 0
$ws.Range("K227").Value = 0
$ws.Range("L227").Value = 0
$ws.Range("M227").Value = 0
$ws.Range("N227").Value = 0
$ws.Range("O227").Value = 0
$ws.Range("P227").Value = 0
$ws.Range("Q227").Value = 0
$ws.Range("R227").Value = 0
$ws.Range("S227").Value = 0
$ws.Range("T227").Value = 0
$ws.Range("U227").Value = 0
$ws.Range("V227").Value = 0
$ws.Range("W227").Value = 0
$ws.Range("X227").Value = 0
$ws.Range("Y227").Value = 0
$ws.Range("Z227").Value = 0
$ws.Range("AA227").Value = 0
$ws.Range("AB227").Value = 1
$ws.Range("AC227").Value = 1
$ws.Range("AD227").Value = 0
$ws.Range("B228").Value = 0
$ws.Range("C228").Value = 0
$ws.Range("D228").Value = 1
$ws.Range("E228").Value = 1
$ws.Range("F228").Value = 1
$ws.Range("G228").Value = 0
$ws.Range("H228").Value = 0
$ws.Range("I228").Value = 0
$ws.Range("J228").Value = 0
$ws.Range("K228").Value = 0
$ws.Range("L228").Value = 0
$ws.Range("M228").Value = 0
$ws.Range("N228").Value = 0
$ws.Range("O228").Value = 0
$ws.Range("P228").Value = 0
$ws.Range("Q228").Value = 0
$ws.Range("R228").Value = 0
$ws.Range("S228").Value = 0
$ws.Range("T228").Value = 0
$ws.Range("U228").Value = 0
$ws.Range("V228").Value = 0
$ws.Range("W228").Value = 0
$ws.Range("X228").Value = 0
$ws.Range("Y228").Value = 0
$ws.Range("Z228").Value = 0
$ws.Range("AA228").Value = 0
$ws.Range("AB228").Value = 1
$ws.Range("AC228").Value = 1
$ws.Range("AD228").Value = 0
$ws.Range("B229").Value = 0
$ws.Range("C229").Value = 0
$ws.Range("D229").Value = 1
$ws.Range("E229").Value = 1
$ws.Range("F229").Value = 1
$ws.Range("G229").Value = 0
$ws.Range("H229").Value = 0
$ws.Range("I229").Value = 0
$ws.Range("J229").Value = 0
$ws.Range("K229").Value = 0
$ws.Range("L229").Value = 0
$ws.Range("M229").Value = 0
$ws.Range("N229").Value = 0
$ws.Range("O229").Value = 0
$ws.Range("P229").Value = 0
$ws.Range("Q229").Value = 0
$ws.Range("R229").Value = 0
$ws.Range("S229").Value = 0
$ws.Range("T229").Value = 0
$ws.Range("U229").Value = 0
$ws.Range("V229").Value = 0
$ws.Range("W229").Value = 0
$ws.Range("X229").Value = 0
$ws.Range("Y229").Value = 0
$ws.Range("Z229").Value = 0
$ws.Range("AA229").Value = 0
$ws.Range("AB229").Value = 1
$ws.Range("AC229").Value = 1
$ws.Range("AD229").Value = 0
$ws.Range("B230").Value = 0
$ws.Range("C230").Value = 0
$ws.Range("D230").Value = 1
$ws.Range("E230").Value = 1
$ws.Range("F230").Value = 1
$ws.Range("G230").Value = 0
$ws.Range("H230").Value = 0
$ws.Range("I230").Value = 0
$ws.Range("J230").Value = 0
$ws.Range("K230").Value = 0
$ws.Range("L230").Value = 0
$ws.Range("M230").Value = 0
$ws.Range("N230").Value = 0
$ws.Range("O230").Value = 0
$ws.Range("P230").Value = 0
$ws.Range("Q230").Value = 0
$ws.Range("R230").Value = 0
$ws.Range("S230").Value = 0
$ws.Range("T230").Value = 0
$ws.Range("U230").Value = 0
$ws.Range("V230").Value = 0
$ws.Range("W230").Value = 0
$ws.Range("X230").Value = 0
$ws.Range("Y230").Value = 0
$ws.Range("Z230").Value = 0
$ws.Range("AA230").Value = 0
$ws.Range("AB230").Value = 1
$ws.Range("AC230").Value = 1
$ws.Range("AD230").Value = 0
$ws.Range("B231").Value = 0
$ws.Range("C231").Value = 0
$ws.Range("D231").Value = 1
$ws.Range("E231").Value = 1
$ws.Range("F231").Value = 1
$ws.Range("G231").Value = 0
$ws.Range("H231").Value = 0
$ws.Range("I231").Value = 0
$ws.Range("J231").Value = 0
$ws.Range("K231").Value = 0
$ws.Range("L231").Value = 0
$ws.Range("M231").Value = 0
$ws.Range("N231").Value = 0
$ws.Range("O231").Value = 0
$ws.Range("P231").Value = 0
$ws.Range("Q231").Value = 0
$ws.Range("R231").Value = 0
$ws.Range("S231").Value = 0
$ws.Range("T231").Value = 0
$ws.Range("U231").Value = 0
$ws.Range("V231").Value = 0
$ws.Range("W231").Value = 0
$ws.Range("X231").Value = 0
$ws.Range("Y231").Value = 0
$ws.Range("Z231").Value = 0
$ws.Range("AA231").Value = 0
$ws.Range("AB231").Value = 1
$ws.Range("AC231").Value = 1
$ws.Range("AD231").Value = 0
$ws.Range("B232").Value = 0
$ws.Range("C232").Value = 0
$ws.Range("D232").Value = 1
$ws.Range("E232").Value = 1
$ws.Range("F232").Value = 1
$ws.Range("G232").Value = 0
$ws.Range("H232").Value = 0
$ws.Range("I232").Value = 0
$ws.Range("J232").Value = 0
$ws.Range("K232").Value = 0
$ws.Range("L232").Value = 0
$ws.Range("M232").Value = 0
$ws.Range("N232").Value = 0
$ws.Range("O232").Value = 0
$ws.Range("P232").Value = 0
$ws.Range("Q232").Value = 0
$ws.Range("R232").Value = 0
$ws.Range("S232").Value = 0
$ws.Range("T232").Value = 0
$ws.Range("U232").Value = 0
$ws.Range("V232").Value = 0
$ws.Range("W232").Value = 0
$ws.Range("X232").Value = 0
$ws.Range("Y232").Value = 0
$ws.Range("Z232").Value = 0
$ws.Range("AA232").Value = 0
$ws.Range("AB232").Value = 1
$ws.Range("AC232").Value = 1
$ws.Range("AD232").Value = 0
$ws.Range("B233").Value = 0
$ws.Range("C233").Value = 0
$ws.Range("D233").Value = 1
$ws.Range("E233").Value = 1
$ws.Range("F233").Value = 1
$ws.Range("G233").Value = 0
$ws.Range("H233").Value = 0
$ws.Range("I233").Value = 0
$ws.Range("J233").Value = 0
$ws.Range("K233").Value = 0
$ws.Range("L233").Value = 0
$ws.Range("M233").Value = 0
$ws.Range("N233").Value = 0
$ws.Range("O233").Value = 0
$ws.Range("P233").Value = 0
$ws.Range("Q233").Value = 0
$ws.Range("R233").Value = 0
$ws.Range("S233").Value = 0
$ws.Range("T233").Value = 0
$ws.Range("U233").Value = 0
$ws.Range("V233").Value = 0
$ws.Range("W233").Value = 0
$ws.Range("X233").Value = 0
$ws.Range("Y233").Value = 0
$ws.Range("Z233").Value = 0
$ws.Range("AA233").Value = 0
$ws.Range("AB233").Value = 1
$ws.Range("AC233").Value = 1
$ws.Range("AD233").Value = 0

# Step 3c: LockdownEffectiveness (AE) for the new rows
$ws.Range("AE222").Value = 0.277777777775
$ws.Range("AE223").Value = 0.277777777775
$ws.Range("AE224").Value = 0.277777777775
$ws.Range("AE225").Value = 0.277777777775
$ws.Range("AE226").Value = 0.277777777775
$ws.Range("AE227").Value = 0.277777777775
$ws.Range("AE228").Value = 0.277777777775
$ws.Range("AE229").Value = 0.277777777775
$ws.Range("AE230").Value = 0.277777777775
$ws.Range("AE231").Value = 0.277777777775
$ws.Range("AE232").Value = 0.277777777775
$ws.Range("AE233").Value = 0.277777777775
